$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 39
$ws.Range("H39").Value = 108.375
$ws.Range("I39").Value = 88.14286
$ws.Range("K39").Value = 264.42858
$ws.Range("M39").Value = 31.57141999999999
# row 86
$ws.Range("H86").Value = 1340.5
$ws.Range("I86").Value = 1165.1111
$ws.Range("K86").Value = 1165.1111
$ws.Range("M86").Value = -42.11110000000008
# row 89
$ws.Range("H89").Value = 1340.5
$ws.Range("I89").Value = 1165.1111
$ws.Range("K89").Value = 5825.5555
$ws.Range("M89").Value = -209.5555000000004
# row 98
$ws.Range("H98").Value = 700397.9399999999
$ws.Range("I98").Value = 1118255.4
$ws.Range("K98").Value = 1118255.4
$ws.Range("M98").Value = -1116757.4
# row 112
$ws.Range("H112").Value = 16043808
$ws.Range("J112").Value = 16043808
$ws.Range("L112").Value = 48131424
$ws.Range("N112").Value = -48133640
# row 113
$ws.Range("H113").Value = 102489.4
$ws.Range("I113").Value = 145227.86
$ws.Range("J113").Value = 2766.3333
$ws.Range("K113").Value = 145227.86
$ws.Range("L113").Value = 2766.3333
$ws.Range("M113").Value = -141973.86
$ws.Range("N113").Value = -9274.3333
# row 122
$ws.Range("H122").Value = 700397.9399999999
$ws.Range("I122").Value = 1118255.4
$ws.Range("K122").Value = 3354766.2
$ws.Range("M122").Value = -3352316.2
# row 129
$ws.Range("H129").Value = 896.6579
$ws.Range("I129").Value = 415.07693
$ws.Range("J129").Value = 1147.08
$ws.Range("K129").Value = 1245.23079
$ws.Range("L129").Value = 3441.24
$ws.Range("M129").Value = 3754.76921
$ws.Range("N129").Value = -13441.24
# row 138
$ws.Range("H138").Value = 5982651.5
$ws.Range("I138").Value = 2151571.8
$ws.Range("J138").Value = 6759221.5
$ws.Range("K138").Value = 6454715.399999999
$ws.Range("L138").Value = 20277664.5
$ws.Range("M138").Value = -6449575.399999999
$ws.Range("N138").Value = -20287944.5

$ws = $wb.Worksheets.Item("ARM")
# row 12
$ws.Range("H12").Value = 50000
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()
# row 17
$ws.Range("H17").Value = 4000
$ws.Range("J17").Value = 4000
$ws.Range("L17").Value = 4000
$ws.Range("N17").Value = -4346
# row 32
$ws.Range("H32").Value = 3147.44
$ws.Range("I32").Value = 1860.439
$ws.Range("J32").Value = 9010.444
$ws.Range("K32").Value = 1860.439
$ws.Range("L32").Value = 9010.444
$ws.Range("M32").Value = -1573.439
$ws.Range("N32").Value = -9584.444
# row 61
$ws.Range("H61").Value = 2687.7036
$ws.Range("I61").Value = 2038.8096
$ws.Range("K61").Value = 2038.8096
$ws.Range("M61").Value = -1826.8096
# row 136
$ws.Range("H136").Value = 2687.7036
$ws.Range("I136").Value = 2038.8096
$ws.Range("K136").Value = 6116.4288
$ws.Range("M136").Value = -3566.4288

$ws = $wb.Worksheets.Item("BSM")
# row 107
$ws.Range("H107").Value = 575.5714
$ws.Range("I107").Value = 544.3333
$ws.Range("J107").Value = 617.2222
$ws.Range("K107").Value = 544.3333
$ws.Range("L107").Value = 617.2222
$ws.Range("M107").Value = 1375.6667
$ws.Range("N107").Value = -4457.2222

$ws = $wb.Worksheets.Item("CRP")
# row 122
$ws.Range("H122").Value = 1412.4
$ws.Range("I122").Value = 1412.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4237.200000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1787.200000000001
$ws.Range("N122").ClearContents()
# row 132
$ws.Range("H132").Value = 3244.5557
$ws.Range("I132").Value = 2705.5715
$ws.Range("J132").Value = 3587.5454
$ws.Range("K132").Value = 8116.7145
$ws.Range("L132").Value = 10762.6362
$ws.Range("M132").Value = -5586.7145
$ws.Range("N132").Value = -15822.6362

$ws = $wb.Worksheets.Item("CUL")
# row 10
$ws.Range("H10").Value = 301
$ws.Range("I10").Value = 301
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 903
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -764
$ws.Range("N10").ClearContents()
# row 17
$ws.Range("H17").Value = 700
$ws.Range("J17").Value = 700
$ws.Range("L17").Value = 2100
$ws.Range("N17").Value = -2438
# row 60
$ws.Range("H60").Value = 3709.5
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 3709.5
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 11128.5
$ws.Range("N60").Value = -11630.5
$ws.Range("M60").ClearContents()
# row 113
$ws.Range("H113").Value = 13514313
$ws.Range("I113").Value = 575.75
$ws.Range("J113").Value = 15152342
$ws.Range("K113").Value = 1727.25
$ws.Range("L113").Value = 45457026
$ws.Range("M113").Value = 442.75
$ws.Range("N113").Value = -45461366
# row 121
$ws.Range("H121").Value = 290
$ws.Range("I121").Value = 290
$ws.Range("K121").Value = 870
$ws.Range("M121").Value = 440

$ws = $wb.Worksheets.Item("GSM")
# row 11
$ws.Range("H11").Value = 11642858
$ws.Range("I11").Value = 11642858
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 11642858
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -11642719
$ws.Range("N11").ClearContents()
# row 80
$ws.Range("H80").Value = 2415.3684
$ws.Range("I80").Value = 2337.7646
$ws.Range("J80").Value = 3075
$ws.Range("K80").Value = 2337.7646
$ws.Range("L80").Value = 3075
$ws.Range("M80").Value = -1339.7646
$ws.Range("N80").Value = -5071
# row 83
$ws.Range("H83").Value = 2415.3684
$ws.Range("I83").Value = 2337.7646
$ws.Range("J83").Value = 3075
$ws.Range("K83").Value = 11688.823
$ws.Range("L83").Value = 15375
$ws.Range("M83").Value = -6696.823
$ws.Range("N83").Value = -25359
# row 102
$ws.Range("H102").Value = 1522.3529
$ws.Range("I102").Value = 1165.0834
$ws.Range("J102").Value = 2379.8
$ws.Range("K102").Value = 1165.0834
$ws.Range("L102").Value = 2379.8
$ws.Range("M102").Value = 456.9166
$ws.Range("N102").Value = -5623.8
# row 107
$ws.Range("H107").Value = 277.54544
$ws.Range("I107").Value = 253
$ws.Range("J107").Value = 298
$ws.Range("K107").Value = 253
$ws.Range("L107").Value = 298
$ws.Range("M107").Value = 1667
$ws.Range("N107").Value = -4138
# row 132
$ws.Range("H132").Value = 2517.8667
$ws.Range("I132").Value = 2173.818
$ws.Range("J132").Value = 3464
$ws.Range("K132").Value = 6521.454000000001
$ws.Range("L132").Value = 10392
$ws.Range("M132").Value = -3991.454000000001
$ws.Range("N132").Value = -15452

$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 11686.667
$ws.Range("I22").Value = 600
$ws.Range("J22").Value = 25545
$ws.Range("K22").Value = 600
$ws.Range("L22").Value = 25545
$ws.Range("M22").Value = -305
$ws.Range("N22").Value = -26135
# row 27
$ws.Range("H27").Value = 11686.667
$ws.Range("I27").Value = 600
$ws.Range("J27").Value = 25545
$ws.Range("K27").Value = 600
$ws.Range("L27").Value = 25545
$ws.Range("M27").Value = -493
$ws.Range("N27").Value = -25759
# row 82
$ws.Range("H82").Value = 1257
$ws.Range("I82").Value = 908.4
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 908.4
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -547.4
$ws.Range("N82").Value = -3722
# row 85
$ws.Range("H85").Value = 1257
$ws.Range("I85").Value = 908.4
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 908.4
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = 339.6
$ws.Range("N85").Value = -5496

$ws = $wb.Worksheets.Item("WVR")
# row 2
$ws.Range("H2").Value = 25013762
$ws.Range("J2").Value = 15471.818
$ws.Range("L2").Value = 15471.818
$ws.Range("N2").Value = -15695.818
# row 122
$ws.Range("H122").Value = 1627.7778
$ws.Range("I122").Value = 2150
$ws.Range("J122").Value = 1366.6666
$ws.Range("K122").Value = 6450
$ws.Range("L122").Value = 4099.9998
$ws.Range("M122").Value = -4000
$ws.Range("N122").Value = -8999.9998
# row 132
$ws.Range("H132").Value = 16671048
$ws.Range("I132").Value = 22731824
$ws.Range("J132").Value = 3910.375
$ws.Range("K132").Value = 68195472
$ws.Range("L132").Value = 11731.125
$ws.Range("M132").Value = -68192942
$ws.Range("N132").Value = -16791.125

